# Generate Report for Handoff
#
# The c4af78f5-b021-4e06-8beb-ab1e89a4191d.md file has finished its
# "In Translation" / "ht" (human-translation) pass and is now
# "Ready for handoff" using machine translation ("mt"), with a refreshed
# handoff datetime. Reflect this on all three sheets:
#   - Overview: zh-cn / de-de status columns + the shared "Latest HO Xliff
#     Generate Date" column for that row.
#   - zh-cn / de-de detail sheets: Status, Priority and Latest Handoff
#     Datetime columns for that row.

$wb = $excel.ActiveWorkbook

$newStatus   = "Ready for handoff"
$newPriority = "mt"
$zhDatetime  = "2016-08-24 02:14:15"
$deDatetime  = "2016-08-24 02:14:19"

# ---- Overview sheet (row 3 = c4af78f5-...) ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus
$ov.Range("G3").Value = $deDatetime

# ---- zh-cn sheet (row 3 = c4af78f5-...) ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = $newStatus
$zh.Range("E3").Value = $newPriority
$zh.Range("H3").Value = $zhDatetime

# ---- de-de sheet (row 3 = c4af78f5-...) ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = $newStatus
$de.Range("E3").Value = $newPriority
$de.Range("H3").Value = $deDatetime

# ---- Column widths widen to fit "Ready for handoff" ----
$ov.Columns.Item(5).ColumnWidth = 16.4
$ov.Columns.Item(6).ColumnWidth = 16.4
$zh.Columns.Item(3).ColumnWidth = 16.4
$de.Columns.Item(3).ColumnWidth = 16.4
